$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, pushing existing rows 5-46 down to 6-47.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 'Vega Modelo de Temuco'
$ws.Range("C5").Value = 'La Araucanía'
$ws.Range("D5").Value = 44750
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112010
$ws.Range("G5").Value = 'Achicoria'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 55
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("N5").Value = '$/caja 18 unidades'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 444
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = 'Hortaliza'
